$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "91.163.12"
$ws.Range("E2").Value = "  +4.25%  "

Set-TextValue $ws.Range("D3") "3.190.78"
$ws.Range("E3").Value = "  +1.00%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.13%  "

Set-TextValue $ws.Range("D5") "218.87"
$ws.Range("E5").Value = "  +5.39%  "

Set-TextValue $ws.Range("D6") "634.35"
$ws.Range("E6").Value = "  +4.57%  "

Set-TextValue $ws.Range("D7") "0.403"
$ws.Range("E7").Value = "  +3.53%  "

Set-TextValue $ws.Range("D8") "0.711"
$ws.Range("E8").Value = "  +5.42%  "

Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  +0.09%  "

Set-TextValue $ws.Range("D10") "3.190.15"
$ws.Range("E10").Value = "  +1.28%  "

Set-TextValue $ws.Range("D11") "0.571"
$ws.Range("E11").Value = "  +6.47%  "

Set-TextValue $ws.Range("D12") "0.181"
$ws.Range("E12").Value = "  +3.31%  "

Set-TextValue $ws.Range("D13") "0.0000259"
$ws.Range("E13").Value = "  +5.80%  "

Set-TextValue $ws.Range("D14") "5.39"
$ws.Range("E14").Value = "  +2.80%  "

Set-TextValue $ws.Range("D15") "90.787.08"
$ws.Range("E15").Value = "  +4.23%  "

Set-TextValue $ws.Range("D16") "3.777.32"
$ws.Range("E16").Value = "  +1.07%  "

Set-TextValue $ws.Range("D17") "33.03"
$ws.Range("E17").Value = "  +2.95%  "

Set-TextValue $ws.Range("D18") "3.201.46"
$ws.Range("E18").Value = "  +2.12%  "

Set-TextValue $ws.Range("D19") "0.0000225"
$ws.Range("E19").Value = "  +70.35%  "

Set-TextValue $ws.Range("D20") "3.33"
$ws.Range("E20").Value = "  +4.64%  "

Set-TextValue $ws.Range("D21") "439.23"
$ws.Range("E21").Value = "  +6.39%  "

Set-TextValue $ws.Range("D22") "13.38"
$ws.Range("E22").Value = "  +0.00%  "

Set-TextValue $ws.Range("D23") "8.55"
$ws.Range("E23").Value = "  +1.15%  "

Set-TextValue $ws.Range("D24") "5.02"
$ws.Range("E24").Value = "  -0.20%  "

Set-TextValue $ws.Range("D25") "5.27"
$ws.Range("E25").Value = "  +2.16%  "

Set-TextValue $ws.Range("D26") "11.77"
$ws.Range("E26").Value = "  -1.53%  "

Set-TextValue $ws.Range("D27") "80.64"
$ws.Range("E27").Value = "  +10.22%  "

Set-TextValue $ws.Range("D28") "3.364.15"
$ws.Range("E28").Value = "  +1.24%  "

Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.14%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.38%  "

Set-TextValue $ws.Range("D31") "0.157"
$ws.Range("E31").Value = "  -1.76%  "

Set-TextValue $ws.Range("D32") "4.14"
$ws.Range("E32").Value = "  +36.74%  "

Set-TextValue $ws.Range("D33") "8.37"
$ws.Range("E33").Value = "  +2.05%  "

Set-TextValue $ws.Range("D34") "529.15"
$ws.Range("E34").Value = "  -3.00%  "

Set-TextValue $ws.Range("D35") "7.00"
$ws.Range("E35").Value = "  +3.27%  "

Set-TextValue $ws.Range("D36") "1.90"
$ws.Range("E36").Value = "  +2.83%  "

Set-TextValue $ws.Range("D37") "1.29"
$ws.Range("E37").Value = "  -1.19%  "

Set-TextValue $ws.Range("D38") "22.46"
$ws.Range("E38").Value = "  +3.12%  "

Set-TextValue $ws.Range("D39") "22.39"
$ws.Range("E39").Value = "  +2.57%  "

Set-TextValue $ws.Range("D40") "1.00"
$ws.Range("E40").Value = "  +0.49%  "

Set-TextValue $ws.Range("D41") "0.126"
$ws.Range("E41").Value = "  -3.23%  "

Set-TextValue $ws.Range("D42") "1.94"
$ws.Range("E42").Value = "  +1.45%  "

Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  -0.03%  "

Set-TextValue $ws.Range("D44") "0.372"
$ws.Range("E44").Value = "  +1.14%  "

Set-TextValue $ws.Range("D45") "147.36"
$ws.Range("E45").Value = "  -1.56%  "

Set-TextValue $ws.Range("D46") "44.17"
$ws.Range("E46").Value = "  +2.53%  "

Set-TextValue $ws.Range("D47") "172.28"
$ws.Range("E47").Value = "  -0.27%  "

Set-TextValue $ws.Range("D48") "0.126"
$ws.Range("E48").Value = "  +0.82%  "

Set-TextValue $ws.Range("D49") "0.748"
$ws.Range("E49").Value = "  +8.27%  "

Set-TextValue $ws.Range("D50") "24.91"
$ws.Range("E50").Value = "  +5.70%  "

Set-TextValue $ws.Range("D51") "1.23"
$ws.Range("E51").Value = "  +0.35%  "
